$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old "label" heading from A1 (no longer used in the new layout)
$ws.Range("A1").ClearContents()

# Column A list (left list), now re-ordered / re-worded, with gaps for new rows
$ws.Range("A2").Value = "g2q,ED,01,a1"
$ws.Range("A3").Value = "333 ED 02 A2"
$ws.Range("A4").Value = "502 ED 01 a2"
$ws.Range("A5").ClearContents()
$ws.Range("A6").Value = "555 aa 01 a1"
$ws.Range("A7").Value = "456 aa 01 a2"
$ws.Range("A8").Value = "www dd 33 b1"
$ws.Range("A13").Value = "stp dd 33b2"
$ws.Range("A16").Value = "444ED02A1"

# New "todo" / "done" headers in row 1
$ws.Range("D1").Value = "todo"
$ws.Range("G1").Value = "done"

# Column D (todo list)
$ws.Range("D2").Value = "12q ED 01 a1"
$ws.Range("D3").Value = "333 ED 02 A2"
$ws.Range("D4").Value = "12q ED 01 a2"
$ws.Range("D5").Value = "444 ED 02 A1"
$ws.Range("D6").Value = "555 aa 01 a1"
$ws.Range("D7").Value = "456 aa 01 a2"

# Column G (done list)
$ws.Range("G2").Value = "12q ED 01 a1"
$ws.Range("G3").Value = "333 ED 02 A2"
$ws.Range("G4").Value = "12q ED 01 a2"
$ws.Range("G5").Value = "444 ED 02 A1"
$ws.Range("G6").Value = "555 aa 01 a1"
$ws.Range("G7").Value = "456 aa 01 a2"

# Update selection to match the saved view state
$null = $ws.Range("A4").Select()
